$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = '3D_Printer_test_fixed_stl_3rd_gen'
$ws.Cells.Item(2, 2).Value = 48701.611744031121
$ws.Cells.Item(2, 3).Value = 42737.72194506447
$ws.Cells.Item(2, 4).Value = 60.522235870361328
$ws.Cells.Item(2, 5).Value = 19552

$ws.Cells.Item(3, 1).Value = 'acoustic_plate'
$ws.Cells.Item(3, 2).Value = 313750.00000000017
$ws.Cells.Item(3, 3).Value = 74000
$ws.Cells.Item(3, 4).Value = 25
$ws.Cells.Item(3, 5).Value = 62430

$ws.Cells.Item(4, 1).Value = 'Another_dremel_bit_holder'
$ws.Cells.Item(4, 2).Value = 143591.84771055059
$ws.Cells.Item(4, 3).Value = 50565.076410959482
$ws.Cells.Item(4, 4).Value = 15
$ws.Cells.Item(4, 5).Value = 40092

$ws.Cells.Item(5, 1).Value = 'bathroom_paper_roll_holder_with_shelf'
$ws.Cells.Item(5, 2).Value = 632926.45120755211
$ws.Cells.Item(5, 3).Value = 103166.7185681635
$ws.Cells.Item(5, 4).Value = 185.93342590332031
$ws.Cells.Item(5, 5).Value = 148381

$ws.Cells.Item(6, 1).Value = 'Beam_100x30mm'
$ws.Cells.Item(6, 2).Value = 90000
$ws.Cells.Item(6, 3).Value = 13800
$ws.Cells.Item(6, 4).Value = 100
$ws.Cells.Item(6, 5).Value = 20559

$ws.Cells.Item(7, 1).Value = 'Beam_30x100mm'
$ws.Cells.Item(7, 2).Value = 90000
$ws.Cells.Item(7, 3).Value = 13800
$ws.Cells.Item(7, 4).Value = 30
$ws.Cells.Item(7, 5).Value = 18452

$ws.Cells.Item(8, 1).Value = 'Cap'
$ws.Cells.Item(8, 2).Value = 5378.5681143174716
$ws.Cells.Item(8, 3).Value = 2548.046221391171
$ws.Cells.Item(8, 4).Value = 12.50000035762787
$ws.Cells.Item(8, 5).Value = 1883

$ws.Cells.Item(9, 1).Value = 'case'
$ws.Cells.Item(9, 2).Value = 12756.302435755561
$ws.Cells.Item(9, 3).Value = 14777.15722902477
$ws.Cells.Item(9, 4).Value = 21.79999923706055
$ws.Cells.Item(9, 5).Value = 6165

$ws.Cells.Item(10, 1).Value = 'Cube_50mm'
$ws.Cells.Item(10, 2).Value = 125000
$ws.Cells.Item(10, 3).Value = 15000
$ws.Cells.Item(10, 4).Value = 50
$ws.Cells.Item(10, 5).Value = 25012

$ws.Cells.Item(11, 1).Value = 'Cute_axolotl_v5'
$ws.Cells.Item(11, 2).Value = 545347.3514940599
$ws.Cells.Item(11, 3).Value = 64131.482762952903
$ws.Cells.Item(11, 4).Value = 76.063808441162109
$ws.Cells.Item(11, 5).Value = 95169

$ws.Cells.Item(12, 1).Value = 'Cylinder'
$ws.Cells.Item(12, 2).Value = 10004.086661324771
$ws.Cells.Item(12, 3).Value = 10260.452771781769
$ws.Cells.Item(12, 4).Value = 56.5
$ws.Cells.Item(12, 5).Value = 5960

$ws.Cells.Item(13, 1).Value = 'Cylinder_10x50mm'
$ws.Cells.Item(13, 2).Value = 3924.3010299656939
$ws.Cells.Item(13, 3).Value = 1727.4993534451301
$ws.Cells.Item(13, 4).Value = 50
$ws.Cells.Item(13, 5).Value = 3493

$ws.Cells.Item(14, 1).Value = 'Cylinder_50x50mm'
$ws.Cells.Item(14, 2).Value = 98152.154952251323
$ws.Cells.Item(14, 3).Value = 11779.615498434559
$ws.Cells.Item(14, 4).Value = 50
$ws.Cells.Item(14, 5).Value = 19236

$ws.Cells.Item(15, 1).Value = 'Dancing_Happy_Dragon'
$ws.Cells.Item(15, 2).Value = 11907.500194774781
$ws.Cells.Item(15, 3).Value = 8387.9595668518687
$ws.Cells.Item(15, 4).Value = 53.611817426979542
$ws.Cells.Item(15, 5).Value = 6038

$ws.Cells.Item(16, 1).Value = 'Dewalt_Drill_Holder'
$ws.Cells.Item(16, 2).Value = 53360.074326238268
$ws.Cells.Item(16, 3).Value = 17634.267536148891
$ws.Cells.Item(16, 4).Value = 20.914825439453121
$ws.Cells.Item(16, 5).Value = 13314

$ws.Cells.Item(17, 1).Value = 'dragon_statue'
$ws.Cells.Item(17, 2).Value = 142432.9869418719
$ws.Cells.Item(17, 3).Value = 27942.80219278032
$ws.Cells.Item(17, 4).Value = 106.9618110656738
$ws.Cells.Item(17, 5).Value = 40828

$ws.Cells.Item(18, 1).Value = 'Flexi-Rex-improved'
$ws.Cells.Item(18, 2).Value = 22534.426690201821
$ws.Cells.Item(18, 3).Value = 13521.03359153675
$ws.Cells.Item(18, 4).Value = 13.00000619888306
$ws.Cells.Item(18, 5).Value = 9471

$ws.Cells.Item(19, 1).Value = 'Flexy_PLA'
$ws.Cells.Item(19, 2).Value = 36019.599144551321
$ws.Cells.Item(19, 3).Value = 101899.9311060513
$ws.Cells.Item(19, 4).Value = 15
$ws.Cells.Item(19, 5).Value = 26859

$ws.Cells.Item(20, 1).Value = 'full_flexi_possuml'
$ws.Cells.Item(20, 2).Value = 44931.394216805042
$ws.Cells.Item(20, 3).Value = 15658.98141698869
$ws.Cells.Item(20, 4).Value = 43.694000244140618
$ws.Cells.Item(20, 5).Value = 16786

$ws.Cells.Item(21, 1).Value = 'FU_COIN_FINAL_12_31'
$ws.Cells.Item(21, 2).Value = 2065.2194960153552
$ws.Cells.Item(21, 3).Value = 3165.1494227013841
$ws.Cells.Item(21, 4).Value = 3.555999994277955
$ws.Cells.Item(21, 5).Value = 944

$ws.Cells.Item(22, 1).Value = 'HairTieBobbyPinBox'
$ws.Cells.Item(22, 2).Value = 58307.285501919381
$ws.Cells.Item(22, 3).Value = 42435.882811048607
$ws.Cells.Item(22, 4).Value = 52
$ws.Cells.Item(22, 5).Value = 26629

$ws.Cells.Item(23, 1).Value = 'Low_poly_otter_3'
$ws.Cells.Item(23, 2).Value = 6129.6016807130982
$ws.Cells.Item(23, 3).Value = 3664.1685254104809
$ws.Cells.Item(23, 4).Value = 47.136669907300757
$ws.Cells.Item(23, 5).Value = 3901

$ws.Cells.Item(24, 1).Value = 'Modern_Geometric_Planter'
$ws.Cells.Item(24, 2).Value = 81801.968693185336
$ws.Cells.Item(24, 3).Value = 35378.383010706508
$ws.Cells.Item(24, 4).Value = 61.018611907958977
$ws.Cells.Item(24, 5).Value = 29981

$ws.Cells.Item(25, 1).Value = 'name_plate'
$ws.Cells.Item(25, 2).Value = 68844.202533932723
$ws.Cells.Item(25, 3).Value = 40335.801494839157
$ws.Cells.Item(25, 4).Value = 38.881500244140618
$ws.Cells.Item(25, 5).Value = 32891

$ws.Cells.Item(26, 1).Value = 'Octopus_with_top_hat_for_some_reason'
$ws.Cells.Item(26, 2).Value = 19580.99561413093
$ws.Cells.Item(26, 3).Value = 15701.678312169301
$ws.Cells.Item(26, 4).Value = 35
$ws.Cells.Item(26, 5).Value = 11148

$ws.Cells.Item(27, 1).Value = 'PCIE_Display_Stand4'
$ws.Cells.Item(27, 2).Value = 169890.51320796809
$ws.Cells.Item(27, 3).Value = 66414.046891629958
$ws.Cells.Item(27, 4).Value = 134.96919822692871
$ws.Cells.Item(27, 5).Value = 47868

$ws.Cells.Item(28, 1).Value = 'Philips_OneBlade_Razor_Station'
$ws.Cells.Item(28, 2).Value = 59047.747089468918
$ws.Cells.Item(28, 3).Value = 14156.128298277659
$ws.Cells.Item(28, 4).Value = 21.496932983398441
$ws.Cells.Item(28, 5).Value = 15679

$ws.Cells.Item(29, 1).Value = 'Plate_10mm'
$ws.Cells.Item(29, 2).Value = 399999.99999999988
$ws.Cells.Item(29, 3).Value = 88000
$ws.Cells.Item(29, 4).Value = 10
$ws.Cells.Item(29, 5).Value = 71282

$ws.Cells.Item(30, 1).Value = 'Pyramid_40mm'
$ws.Cells.Item(30, 2).Value = 18475.209554036319
$ws.Cells.Item(30, 3).Value = 4800.0001221054636
$ws.Cells.Item(30, 4).Value = 34.641017913818359
$ws.Cells.Item(30, 5).Value = 4850

$ws.Cells.Item(31, 1).Value = 'RubberDuck'
$ws.Cells.Item(31, 2).Value = 357445.95730624709
$ws.Cells.Item(31, 3).Value = 44629.031865117176
$ws.Cells.Item(31, 4).Value = 87.12722110748291
$ws.Cells.Item(31, 5).Value = 52696

$ws.Cells.Item(32, 1).Value = 'scrapper_v5'
$ws.Cells.Item(32, 2).Value = 33840.568813167083
$ws.Cells.Item(32, 3).Value = 11788.292434716859
$ws.Cells.Item(32, 4).Value = 10
$ws.Cells.Item(32, 5).Value = 8233

$ws.Cells.Item(33, 1).Value = 'Shark_popcorn_bowl'
$ws.Cells.Item(33, 2).Value = 1243998.1307313109
$ws.Cells.Item(33, 3).Value = 125356.8620921487
$ws.Cells.Item(33, 4).Value = 166.05999755859381
$ws.Cells.Item(33, 5).Value = 252230

$ws.Cells.Item(34, 1).Value = 'SM_Buddha_print'
$ws.Cells.Item(34, 2).Value = 297302.65980893082
$ws.Cells.Item(34, 3).Value = 41099.007819179933
$ws.Cells.Item(34, 4).Value = 115.76584625244141
$ws.Cells.Item(34, 5).Value = 72176

$ws.Cells.Item(35, 1).Value = 'Speed_Square'
$ws.Cells.Item(35, 2).Value = 119995.2882208134
$ws.Cells.Item(35, 3).Value = 44779.039457428531
$ws.Cells.Item(35, 4).Value = 185.88043212890619
$ws.Cells.Item(35, 5).Value = 36639

$ws.Cells.Item(36, 1).Value = 'Swan'
$ws.Cells.Item(36, 2).Value = 46790.190569169717
$ws.Cells.Item(36, 3).Value = 10277.512151710391
$ws.Cells.Item(36, 4).Value = 120
$ws.Cells.Item(36, 5).Value = 14371

$ws.Cells.Item(37, 1).Value = 'top'
$ws.Cells.Item(37, 2).Value = 11173.03485777296
$ws.Cells.Item(37, 3).Value = 7161.7369636316753
$ws.Cells.Item(37, 4).Value = 3.8000011444091801
$ws.Cells.Item(37, 5).Value = 2708

$ws.Cells.Item(38, 1).Value = 'VaseV846'
$ws.Cells.Item(38, 2).Value = 468170.01108463207
$ws.Cells.Item(38, 3).Value = 51083.140346867993
$ws.Cells.Item(38, 4).Value = 170
$ws.Cells.Item(38, 5).Value = 102112

$excel.ActiveWindow.Zoom = 100
$ws.Range("I15").Select()